$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = '[Harrison X.%Bai%NULL%2,             Ben%Hsieh%NULL%1,             Zeng%Xiong%NULL%1,             Kasey%Halsey%NULL%1,             Ji Whae%Choi%NULL%1,             Thi My Linh%Tran%NULL%1,             Ian%Pan%NULL%1,             Lin-Bo%Shi%NULL%1,             Dong-Cui%Wang%NULL%1,             Ji%Mei%NULL%1,             Xiao-Long%Jiang%NULL%1,             Qiu-Hua%Zeng%NULL%1,             Thomas K.%Egglin%NULL%1,             Ping-Feng%Hu%NULL%1,             Saurabh%Agarwal%NULL%1,             Fangfang%Xie%NULL%1,             Sha%Li%NULL%1,             Terrance%Healey%NULL%1,             Michael K.%Atalay%NULL%1,             Wei-Hua%Liao%liaoweihua2017@163.com%1,            Harrison X.%Bai%NULL%1,            Ben%Hsieh%NULL%1,            Zeng%Xiong%NULL%1,            Kasey%Halsey%NULL%1,            Ji Whae%Choi%NULL%1,            Thi My Linh%Tran%NULL%1,            Ian%Pan%NULL%1,            Lin-Bo%Shi%NULL%1,            Dong-Cui%Wang%NULL%1,            Ji%Mei%NULL%1,            Xiao-Long%Jiang%NULL%1,            Qiu-Hua%Zeng%NULL%1,            Thomas K.%Egglin%NULL%1,            Ping-Feng%Hu%NULL%1,            Saurabh%Agarwal%NULL%1,            Fangfang%Xie%NULL%1,            Sha%Li%NULL%1,            Terrance%Healey%NULL%1,            Michael K.%Atalay%NULL%1,            Wei-Hua%Liao%liaoweihua2017@163.com%1]'
$ws.Range("E4").Value = '[Damiano%Caruso%NULL%2,             Marta%Zerunian%NULL%1,             Michela%Polici%NULL%1,             Francesco%Pucciarelli%NULL%1,             Tiziano%Polidori%NULL%1,             Carlotta%Rucci%NULL%1,             Gisella%Guido%NULL%1,             Benedetta%Bracci%NULL%1,             Chiara%de Dominicis%NULL%1,             Andrea%Laghi%andrea.laghi@uniroma1.it%1,            Damiano%Caruso%NULL%1,            Marta%Zerunian%NULL%1,            Michela%Polici%NULL%1,            Francesco%Pucciarelli%NULL%1,            Tiziano%Polidori%NULL%1,            Carlotta%Rucci%NULL%1,            Gisella%Guido%NULL%1,            Benedetta%Bracci%NULL%1,            Chiara%de Dominicis%NULL%1,            Andrea%Laghi%andrea.laghi@uniroma1.it%1]'
$ws.Range("E5").Value = '[Xiaofeng%Chen%NULL%2,             Yanyan%Tang%NULL%1,             Yongkang%Mo%NULL%1,             Shengkai%Li%NULL%1,             Daiying%Lin%NULL%1,             Zhijian%Yang%NULL%1,             Zhiqi%Yang%NULL%1,             Hongfu%Sun%NULL%1,             Jinming%Qiu%NULL%1,             Yuting%Liao%NULL%1,             Jianning%Xiao%NULL%1,             Xiangguang%Chen%NULL%1,             Xianheng%Wu%NULL%1,             Renhua%Wu%NULL%1,             Zhuozhi%Dai%zhuozhi@ualberta.ca%1,            Xiaofeng%Chen%NULL%1,            Yanyan%Tang%NULL%1,            Yongkang%Mo%NULL%1,            Shengkai%Li%NULL%1,            Daiying%Lin%NULL%1,            Zhijian%Yang%NULL%1,            Zhiqi%Yang%NULL%1,            Hongfu%Sun%NULL%1,            Jinming%Qiu%NULL%1,            Yuting%Liao%NULL%1,            Jianning%Xiao%NULL%1,            Xiangguang%Chen%NULL%1,            Xianheng%Wu%NULL%1,            Renhua%Wu%NULL%1,            Zhuozhi%Dai%zhuozhi@ualberta.ca%1]'
$ws.Range("E7").Value = '[Hyewon%Choi%NULL%2,             Xiaolong%Qi%NULL%1,             Soon Ho%Yoon%yshoka@gmail.com%1,             Sang Joon%Park%NULL%1,             Kyung Hee%Lee%NULL%1,             Jin Yong%Kim%NULL%1,             Young Kyung%Lee%NULL%1,             Hongseok%Ko%NULL%1,             Ki Hwan%Kim%NULL%1,             Chang Min%Park%NULL%1,             Yun-Hyeon%Kim%NULL%1,             Junqiang%Lei%NULL%1,             Jung Hee%Hong%NULL%1,             Hyungjin%Kim%NULL%1,             Eui Jin%Hwang%NULL%1,             Seung Jin%Yoo%NULL%1,             Ju Gang%Nam%NULL%1,             Chang Hyun%Lee%NULL%1,             Jin Mo%Goo%NULL%1,            Hyewon%Choi%NULL%1,            Xiaolong%Qi%NULL%1,            Soon Ho%Yoon%yshoka@gmail.com%1,            Sang Joon%Park%NULL%1,            Kyung Hee%Lee%NULL%1,            Jin Yong%Kim%NULL%1,            Young Kyung%Lee%NULL%1,            Hongseok%Ko%NULL%1,            Ki Hwan%Kim%NULL%1,            Chang Min%Park%NULL%1,            Yun-Hyeon%Kim%NULL%1,            Junqiang%Lei%NULL%1,            Jung Hee%Hong%NULL%1,            Hyungjin%Kim%NULL%1,            Eui Jin%Hwang%NULL%1,            Seung Jin%Yoo%NULL%1,            Ju Gang%Nam%NULL%1,            Chang Hyun%Lee%NULL%1,            Jin Mo%Goo%NULL%1]'
$ws.Range("E8").Value = '[Yuki%Himoto%yukihimoto@gmail.com%2,             Akihiko%Sakata%NULL%1,             Mitsuhiro%Kirita%NULL%1,             Takashi%Hiroi%NULL%1,             Ken-ichiro%Kobayashi%NULL%1,             Kenji%Kubo%NULL%1,             Hyunjin%Kim%NULL%1,             Azusa%Nishimoto%NULL%1,             Chikara%Maeda%NULL%1,             Akira%Kawamura%NULL%1,             Nobuhiro%Komiya%NULL%1,             Shigeaki%Umeoka%NULL%1,            Yuki%Himoto%yukihimoto@gmail.com%1,            Akihiko%Sakata%NULL%1,            Mitsuhiro%Kirita%NULL%1,            Takashi%Hiroi%NULL%1,            Ken-ichiro%Kobayashi%NULL%1,            Kenji%Kubo%NULL%1,            Hyunjin%Kim%NULL%1,            Azusa%Nishimoto%NULL%1,            Chikara%Maeda%NULL%1,            Akira%Kawamura%NULL%1,            Nobuhiro%Komiya%NULL%1,            Shigeaki%Umeoka%NULL%1]'
$ws.Range("E9").Value = '[Chunqin%Long%NULL%2,             Huaxiang%Xu%NULL%1,             Qinglin%Shen%NULL%1,             Xianghai%Zhang%NULL%1,             Bing%Fan%26171381@qq.com%1,             Chuanhong%Wang%NULL%1,             Bingliang%Zeng%NULL%1,             Zicong%Li%NULL%1,             Xiaofen%Li%NULL%1,             Honglu%Li%NULL%1,            Chunqin%Long%NULL%1,            Huaxiang%Xu%NULL%1,            Qinglin%Shen%NULL%1,            Xianghai%Zhang%NULL%1,            Bing%Fan%26171381@qq.com%1,            Chuanhong%Wang%NULL%1,            Bingliang%Zeng%NULL%1,            Zicong%Li%NULL%1,            Xiaofen%Li%NULL%1,            Honglu%Li%NULL%1]'
$ws.Range("E10").Value = '[Anthony%Dangis%NULL%4,            Christopher%Gieraerts%NULL%1,            Yves%De Bruecker%NULL%1,            Lode%Janssen%NULL%1,            Hanne%Valgaeren%NULL%1,            Dagmar%Obbels%NULL%1,            Marc%Gillis%NULL%1,            Marc%Van Ranst%NULL%0,            Johan%Frans%NULL%1,            Annick%Demeyere%NULL%1,            Rolf%Symons%rolf.symons@imelda.be%1]'
$ws.Range("E11").Value = '[Wanbo%Zhu%NULL%0,          Kai%Xie%NULL%0,          Kai%Xie%NULL%0,          Hui%Lu%NULL%0,          Lei%Xu%bayinhexl@126.com%0,          Shusheng%Zhou%zhouss108@163.com%0,          Shiyuan%Fang%fangshiyuan2008@126.com%0]'
$ws.Range("E12").Value = '[Anthony%Dangis%NULL%0,          Christopher%Gieraerts%NULL%2,          Christopher%Gieraerts%NULL%0,          Yves%De Bruecker%NULL%2,          Yves%De Bruecker%NULL%0,          Lode%Janssen%NULL%2,          Lode%Janssen%NULL%0,          Hanne%Valgaeren%NULL%2,          Hanne%Valgaeren%NULL%0,          Dagmar%Obbels%NULL%2,          Dagmar%Obbels%NULL%0,          Marc%Gillis%NULL%2,          Marc%Gillis%NULL%0,          Marc%Van Ranst%NULL%0,          Marc%Van Ranst%NULL%0,          Johan%Frans%NULL%1,          Annick%Demeyere%NULL%2,          Annick%Demeyere%NULL%0,          Rolf%Symons%rolf.symons@imelda.be%2,          Rolf%Symons%rolf.symons@imelda.be%0]'
$ws.Range("E13").Value = '[Anthony%Dangis%NULL%0,       Christopher%Gieraerts%NULL%2,       Christopher%Gieraerts%NULL%0,       Yves%De Bruecker%NULL%2,       Yves%De Bruecker%NULL%0,       Lode%Janssen%NULL%2,       Lode%Janssen%NULL%0,       Hanne%Valgaeren%NULL%2,       Hanne%Valgaeren%NULL%0,       Dagmar%Obbels%NULL%2,       Dagmar%Obbels%NULL%0,       Marc%Gillis%NULL%2,       Marc%Gillis%NULL%0,       Marc%Van Ranst%NULL%2,       Marc%Van Ranst%NULL%0,       Johan%Frans%NULL%1,       Annick%Demeyere%NULL%2,       Annick%Demeyere%NULL%0,       Rolf%Symons%rolf.symons@imelda.be%2,       Rolf%Symons%rolf.symons@imelda.be%0]'
